$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 4920.722
$ws.Range("I38").Value = 138.5
$ws.Range("J38").Value = 6287.0713
$ws.Range("K38").Value = 415.5
$ws.Range("L38").Value = 18861.2139
$ws.Range("M38").Value = -43.5
$ws.Range("N38").Value = -19605.2139
$ws.Range("H96").Value = 43480756
$ws.Range("I96").Value = 3627.7144
$ws.Range("J96").Value = 111111850
$ws.Range("K96").Value = 10883.1432
$ws.Range("L96").Value = 333335550
$ws.Range("M96").Value = -9510.143199999999
$ws.Range("N96").Value = -333338296
$ws.Range("H112").Value = 2608.0833
$ws.Range("I112").Value = 1738.8
$ws.Range("J112").Value = 3229
$ws.Range("K112").Value = 5216.4
$ws.Range("L112").Value = 9687
$ws.Range("M112").Value = -4108.4
$ws.Range("N112").Value = -11903
$ws.Range("H118").Value = 1173.6364
$ws.Range("I118").Value = 1173.6364
$ws.Range("K118").Value = 3520.9092
$ws.Range("M118").Value = -1863.9092
$ws.Range("H132").Value = 6912.96
$ws.Range("I132").Value = 4366.2285
$ws.Range("K132").Value = 13098.6855
$ws.Range("M132").Value = -10568.6855
$ws.Range("H135").Value = 1158.6666
$ws.Range("I135").Value = 589.3333
$ws.Range("J135").Value = 2866.6667
$ws.Range("K135").Value = 5303.9997
$ws.Range("L135").Value = 25800.0003
$ws.Range("M135").Value = -2768.9997
$ws.Range("N135").Value = -30870.0003
$ws.Range("H137").Value = 3824.262
$ws.Range("I137").Value = 4611.3667
$ws.Range("J137").Value = 1856.5
$ws.Range("K137").Value = 13834.1001
$ws.Range("L137").Value = 5569.5
$ws.Range("M137").Value = -11284.1001
$ws.Range("N137").Value = -10669.5
$ws.Range("H138").Value = 12503826
$ws.Range("I138").Value = 35717544
$ws.Range("J138").Value = 4132.4424
$ws.Range("K138").Value = 107152632
$ws.Range("L138").Value = 12397.3272
$ws.Range("M138").Value = -107147492
$ws.Range("N138").Value = -22677.3272
$ws.Range("H141").Value = 11381.897
$ws.Range("I141").Value = 15095
$ws.Range("J141").Value = 7049.9443
$ws.Range("K141").Value = 45285
$ws.Range("L141").Value = 21149.8329
$ws.Range("M141").Value = -40105
$ws.Range("N141").Value = -31509.8329

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 11905711
$ws.Range("I2").Value = 14286519
$ws.Range("K2").Value = 14286519
$ws.Range("M2").Value = -14286406
$ws.Range("H45").Value = 12325.053
$ws.Range("I45").Value = 11451.077
$ws.Range("K45").Value = 11451.077
$ws.Range("M45").Value = -11074.077
$ws.Range("H61").Value = 6660.0967
$ws.Range("I61").Value = 6355.381
$ws.Range("K61").Value = 6355.381
$ws.Range("M61").Value = -6143.381
$ws.Range("H63").Value = 3083.037
$ws.Range("I63").Value = 2187.0833
$ws.Range("K63").Value = 2187.0833
$ws.Range("M63").Value = -1501.0833
$ws.Range("H66").Value = 3083.037
$ws.Range("I66").Value = 2187.0833
$ws.Range("K66").Value = 10935.4165
$ws.Range("M66").Value = -7503.416499999999
$ws.Range("H88").Value = 1399.909
$ws.Range("I88").Value = 1357.1428
$ws.Range("J88").Value = 1474.75
$ws.Range("K88").Value = 1357.1428
$ws.Range("L88").Value = 1474.75
$ws.Range("M88").Value = -951.1428000000001
$ws.Range("N88").Value = -2286.75
$ws.Range("H91").Value = 1399.909
$ws.Range("I91").Value = 1357.1428
$ws.Range("J91").Value = 1474.75
$ws.Range("K91").Value = 1357.1428
$ws.Range("L91").Value = 1474.75
$ws.Range("M91").Value = 46.85719999999992
$ws.Range("N91").Value = -4282.75
$ws.Range("H116").Value = 11905711
$ws.Range("I116").Value = 14286519
$ws.Range("K116").Value = 14286519
$ws.Range("M116").Value = -14284225
$ws.Range("H136").Value = 6660.0967
$ws.Range("I136").Value = 6355.381
$ws.Range("K136").Value = 19066.143
$ws.Range("M136").Value = -16516.143

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 11905711
$ws.Range("I3").Value = 14286519
$ws.Range("K3").Value = 14286519
$ws.Range("M3").Value = -14286405
$ws.Range("H134").Value = 2611.2778
$ws.Range("I134").Value = 2039.8937
$ws.Range("K134").Value = 6119.6811
$ws.Range("M134").Value = -3584.6811

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3997.087
$ws.Range("I31").Value = 2015.2
$ws.Range("K31").Value = 2015.2
$ws.Range("M31").Value = -1720.2
$ws.Range("H34").Value = 3997.087
$ws.Range("I34").Value = 2015.2
$ws.Range("K34").Value = 2015.2
$ws.Range("M34").Value = -1813.2
$ws.Range("H74").Value = 92499.5
$ws.Range("J74").Value = 92499.5
$ws.Range("L74").Value = 92499.5
$ws.Range("N74").Value = -94247.5
$ws.Range("H77").Value = 92499.5
$ws.Range("J77").Value = 92499.5
$ws.Range("L77").Value = 277498.5
$ws.Range("N77").Value = -286234.5
$ws.Range("H86").Value = 2904156.8
$ws.Range("I86").Value = 4169409.5
$ws.Range("J86").Value = 12150.571
$ws.Range("K86").Value = 4169409.5
$ws.Range("L86").Value = 12150.571
$ws.Range("M86").Value = -4168286.5
$ws.Range("N86").Value = -14396.571
$ws.Range("H89").Value = 2904156.8
$ws.Range("I89").Value = 4169409.5
$ws.Range("J89").Value = 12150.571
$ws.Range("K89").Value = 20847047.5
$ws.Range("L89").Value = 60752.855
$ws.Range("M89").Value = -20841431.5
$ws.Range("N89").Value = -71984.855
$ws.Range("H107").Value = 1652.5883
$ws.Range("I107").Value = 1261.8
$ws.Range("K107").Value = 1261.8
$ws.Range("M107").Value = 658.2
$ws.Range("H109").Value = 144142
$ws.Range("J109").Value = 144142
$ws.Range("L109").Value = 144142
$ws.Range("N109").Value = -146222
$ws.Range("H132").Value = 2780.7144
$ws.Range("I132").Value = 2137.0625
$ws.Range("J132").Value = 4840.4
$ws.Range("K132").Value = 6411.1875
$ws.Range("L132").Value = 14521.2
$ws.Range("M132").Value = -3881.1875
$ws.Range("N132").Value = -19581.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H74").Value = 3859.3333
$ws.Range("I74").Value = 3859.3333
$ws.Range("K74").Value = 11577.9999
$ws.Range("M74").Value = -10516.9999
$ws.Range("H77").Value = 3859.3333
$ws.Range("I77").Value = 3859.3333
$ws.Range("K77").Value = 34733.9997
$ws.Range("M77").Value = -29429.9997
$ws.Range("H136").Value = 2249.6667
$ws.Range("I136").Value = 2249.6667
$ws.Range("K136").Value = 6749.000100000001
$ws.Range("M136").Value = -1649.000100000001
$ws.Range("H137").Value = 5049.3076
$ws.Range("I137").Value = 3649.6
$ws.Range("K137").Value = 10948.8
$ws.Range("M137").Value = -5848.799999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3953.4211
$ws.Range("J61").Value = 1500
$ws.Range("L61").Value = 1500
$ws.Range("N61").Value = -1904
$ws.Range("H113").Value = 3953.4211
$ws.Range("J113").Value = 1500
$ws.Range("L113").Value = 1500
$ws.Range("N113").Value = -5840
$ws.Range("H136").Value = 6213818
$ws.Range("J136").Value = 5672
$ws.Range("L136").Value = 17016
$ws.Range("N136").Value = -22116

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 179966.5
$ws.Range("J16").Value = 179966.5
$ws.Range("L16").Value = 179966.5
$ws.Range("N16").Value = -180550.5
$ws.Range("H62").Value = 20841204
$ws.Range("I62").Value = 22229790
$ws.Range("K62").Value = 22229790
$ws.Range("M62").Value = -22229166
$ws.Range("H65").Value = 20841204
$ws.Range("I65").Value = 22229790
$ws.Range("K65").Value = 111148950
$ws.Range("M65").Value = -111145830
$ws.Range("H136").Value = 2298.8684
$ws.Range("I136").Value = 2178.1765
$ws.Range("J136").Value = 3324.75
$ws.Range("K136").Value = 6534.529500000001
$ws.Range("L136").Value = 9974.25
$ws.Range("M136").Value = -3984.529500000001
$ws.Range("N136").Value = -15074.25
